$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 29.747261
$ws.Range("H2").Value = 59.494522
$ws.Range("I2").Value = 0.1897261270801772
$ws.Range("J2").Value = 0.1356501617082849
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 21.983953
$ws.Range("N2").Value = 65.951859
$ws.Range("O2").Value = 0.3824831516716194
$ws.Range("P2").Value = 0.3824831516716194
$ws.Range("Q2").Value = 653.9623877027329
$ws.Range("R2").Value = 3923.774326216398
$ws.Range("S2").Value = 0.07256704704007635
$ws.Range("T2").Value = 0.05188390137494964

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 29.747261
$ws.Range("H3").Value = 59.494522
$ws.Range("I3").Value = 0.1897261270801772
$ws.Range("J3").Value = 0.1356501617082849
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 35.37535733333333
$ws.Range("N3").Value = 106.126072
$ws.Range("O3").Value = 0.6154706646417534
$ws.Range("P3").Value = 0.6154706646417534
$ws.Range("Q3").Value = 1052.319987562931
$ws.Range("R3").Value = 6313.919925377583
$ws.Range("S3").Value = 0.1167708655339424
$ws.Range("T3").Value = 0.08348869518535942

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("G4").Value = 29.747261
$ws.Range("H4").Value = 59.494522
$ws.Range("I4").Value = 0.1897261270801772
$ws.Range("J4").Value = 0.1356501617082849
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.1176083333333333
$ws.Range("N4").Value = 0.352825
$ws.Range("O4").Value = 0.002046183686627228
$ws.Range("P4").Value = 0.002046183686627228
$ws.Range("Q4").Value = 3.498525787441666
$ws.Range("R4").Value = 20.99115472465
$ws.Range("S4").Value = 0.000388214506158423
$ws.Range("T4").Value = 0.000277565147975838

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 21.75021766666667
$ws.Range("H5").Value = 65.250653
$ws.Range("I5").Value = 0.1387214964445808
$ws.Range("J5").Value = 0.1487743969271857
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 21.983953
$ws.Range("N5").Value = 65.951859
$ws.Range("O5").Value = 0.3824831516716194
$ws.Range("P5").Value = 0.3824831516716194
$ws.Range("Q5").Value = 478.1557629237697
$ws.Range("R5").Value = 4303.401866313927
$ws.Range("S5").Value = 0.05305863516472661
$ws.Range("T5").Value = 0.0569037002247545

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 21.75021766666667
$ws.Range("H6").Value = 65.250653
$ws.Range("I6").Value = 0.1387214964445808
$ws.Range("J6").Value = 0.1487743969271857
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 35.37535733333333
$ws.Range("N6").Value = 106.126072
$ws.Range("O6").Value = 0.6154706646417534
$ws.Range("P6").Value = 0.6154706646417534
$ws.Range("Q6").Value = 769.4217220361129
$ws.Range("R6").Value = 6924.795498325016
$ws.Range("S6").Value = 0.08537901161684476
$ws.Range("T6").Value = 0.09156627695845104

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 21.75021766666667
$ws.Range("H7").Value = 65.250653
$ws.Range("I7").Value = 0.1387214964445808
$ws.Range("J7").Value = 0.1487743969271857
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.1176083333333333
$ws.Range("N7").Value = 0.352825
$ws.Range("O7").Value = 0.002046183686627228
$ws.Range("P7").Value = 0.002046183686627228
$ws.Range("Q7").Value = 2.558006849413889
$ws.Range("R7").Value = 23.022061644725
$ws.Range("S7").Value = 0.0002838496630094182
$ws.Range("T7").Value = 0.0003044197439802115

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 31.00397566666667
$ws.Range("H8").Value = 93.011927
$ws.Range("I8").Value = 0.1977413728048684
$ws.Range("J8").Value = 0.2120713389099788
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 21.983953
$ws.Range("N8").Value = 65.951859
$ws.Range("O8").Value = 0.3824831516716194
$ws.Range("P8").Value = 0.3824831516716194
$ws.Range("Q8").Value = 681.5899438691437
$ws.Range("R8").Value = 6134.309494822293
$ws.Range("S8").Value = 0.07563274348627873
$ws.Range("T8").Value = 0.08111371408550884

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 31.00397566666667
$ws.Range("H9").Value = 93.011927
$ws.Range("I9").Value = 0.1977413728048684
$ws.Range("J9").Value = 0.2120713389099788
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 35.37535733333333
$ws.Range("N9").Value = 106.126072
$ws.Range("O9").Value = 0.6154706646417534
$ws.Range("P9").Value = 0.6154706646417534
$ws.Range("Q9").Value = 1096.776717962305
$ws.Range("R9").Value = 9870.990461660744
$ws.Range("S9").Value = 0.1217040141473851
$ws.Range("T9").Value = 0.1305236879103912

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 31.00397566666667
$ws.Range("H10").Value = 93.011927
$ws.Range("I10").Value = 0.1977413728048684
$ws.Range("J10").Value = 0.2120713389099788
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.1176083333333333
$ws.Range("N10").Value = 0.352825
$ws.Range("O10").Value = 0.002046183686627228
$ws.Range("P10").Value = 0.002046183686627228
$ws.Range("Q10").Value = 3.646325904863889
$ws.Range("R10").Value = 32.816933143775
$ws.Range("S10").Value = 0.0004046151712045947
$ws.Range("T10").Value = 0.0004339369140787928

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 54.68136733333333
$ws.Range("H11").Value = 164.044102
$ws.Range("I11").Value = 0.3487542616983074
$ws.Range("J11").Value = 0.37402786366769
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 21.983953
$ws.Range("N11").Value = 65.951859
$ws.Range("O11").Value = 0.3824831516716194
$ws.Range("P11").Value = 0.3824831516716194
$ws.Range("Q11").Value = 1202.112609431735
$ws.Range("R11").Value = 10819.01348488562
$ws.Range("S11").Value = 0.1333926291732773
$ws.Range("T11").Value = 0.1430593561086209

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 54.68136733333333
$ws.Range("H12").Value = 164.044102
$ws.Range("I12").Value = 0.3487542616983074
$ws.Range("J12").Value = 0.37402786366769
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 35.37535733333333
$ws.Range("N12").Value = 106.126072
$ws.Range("O12").Value = 0.6154706646417534
$ws.Range("P12").Value = 0.6154706646417534
$ws.Range("Q12").Value = 1934.372908891927
$ws.Range("R12").Value = 17409.35618002735
$ws.Range("S12").Value = 0.2146480172441012
$ws.Range("T12").Value = 0.2302031778460883

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 54.68136733333333
$ws.Range("H13").Value = 164.044102
$ws.Range("I13").Value = 0.3487542616983074
$ws.Range("J13").Value = 0.37402786366769
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.1176083333333333
$ws.Range("N13").Value = 0.352825
$ws.Range("O13").Value = 0.002046183686627228
$ws.Range("P13").Value = 0.002046183686627228
$ws.Range("Q13").Value = 6.430984476461111
$ws.Range("R13").Value = 57.87886028815
$ws.Range("S13").Value = 0.0007136152809287996
$ws.Range("T13").Value = 0.0007653297129808603

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 17.571287
$ws.Range("H14").Value = 52.713861
$ws.Range("I14").Value = 0.1120685440694613
$ws.Range("J14").Value = 0.1201899524281925
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 21.983953
$ws.Range("N14").Value = 65.951859
$ws.Range("O14").Value = 0.3824831516716194
$ws.Range("P14").Value = 0.3824831516716194
$ws.Range("Q14").Value = 386.286347557511
$ws.Range("R14").Value = 3476.577128017599
$ws.Range("S14").Value = 0.04286432993893732
$ws.Range("T14").Value = 0.04597063180399708

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 17.571287
$ws.Range("H15").Value = 52.713861
$ws.Range("I15").Value = 0.1120685440694613
$ws.Range("J15").Value = 0.1201899524281925
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 35.37535733333333
$ws.Range("N15").Value = 106.126072
$ws.Range("O15").Value = 0.6154706646417534
$ws.Range("P15").Value = 0.6154706646417534
$ws.Range("Q15").Value = 621.5905564315548
$ws.Range("R15").Value = 5594.315007883992
$ws.Range("S15").Value = 0.06897490130386497
$ws.Range("T15").Value = 0.07397338990424036

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 17.571287
$ws.Range("H16").Value = 52.713861
$ws.Range("I16").Value = 0.1120685440694613
$ws.Range("J16").Value = 0.1201899524281925
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.1176083333333333
$ws.Range("N16").Value = 0.352825
$ws.Range("O16").Value = 0.002046183686627228
$ws.Range("P16").Value = 0.002046183686627228
$ws.Range("Q16").Value = 2.066529778591667
$ws.Range("R16").Value = 18.598768007325
$ws.Range("S16").Value = 0.0002293128266589962
$ws.Range("T16").Value = 0.0002459307199550701

# Row 17
$ws.Range("E17").Value = 2
$ws.Range("G17").Value = 2.0364265
$ws.Range("H17").Value = 4.072853
$ws.Range("I17").Value = 0.0129881979026049
$ws.Range("J17").Value = 0.009286286358668003
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 21.983953
$ws.Range("N17").Value = 65.951859
$ws.Range("O17").Value = 0.3824831516716194
$ws.Range("P17").Value = 0.3824831516716194
$ws.Range("Q17").Value = 44.7687044639545
$ws.Range("R17").Value = 268.612226783727
$ws.Range("S17").Value = 0.00496776686832304
$ws.Range("T17").Value = 0.003551848073788504

# Row 18
$ws.Range("E18").Value = 2
$ws.Range("G18").Value = 2.0364265
$ws.Range("H18").Value = 4.072853
$ws.Range("I18").Value = 0.0129881979026049
$ws.Range("J18").Value = 0.009286286358668003
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 35.37535733333333
$ws.Range("N18").Value = 106.126072
$ws.Range("O18").Value = 0.6154706646417534
$ws.Range("P18").Value = 0.6154706646417534
$ws.Range("Q18").Value = 72.03931512056934
$ws.Range("R18").Value = 432.235890723416
$ws.Range("S18").Value = 0.007993854795614866
$ws.Range("T18").Value = 0.005715436837223043

# Row 19
$ws.Range("E19").Value = 2
$ws.Range("G19").Value = 2.0364265
$ws.Range("H19").Value = 4.072853
$ws.Range("I19").Value = 0.0129881979026049
$ws.Range("J19").Value = 0.009286286358668003
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 0.1176083333333333
$ws.Range("N19").Value = 0.352825
$ws.Range("O19").Value = 0.002046183686627228
$ws.Range("P19").Value = 0.002046183686627228
$ws.Range("Q19").Value = 0.2395007266208333
$ws.Range("R19").Value = 1.437004359725
$ws.Range("S19").Value = 0.00002657623866699613
$ws.Range("T19").Value = 0.00001900144765645543
